$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column D
$ws.Range("D1").Value = "Seniority"

# Fill column D (rows 2-29) with seniority value 4 for every response row
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 4
}

# Update the active selection as reflected in the saved workbook
$ws.Range("E9").Select()
